$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (2021-12-28, date serial 44558) was added to the
# data table. It belongs right before the current row 14, so push the
# existing rows 14..108 down by one (Excel re-indexes all row-relative
# references automatically) and then populate the freshly inserted row.
$ws.Rows.Item(14).Insert()

$ws.Range("A14").Value = 8
$ws.Range("B14").Value = "Terminal La Palmera de La Serena"
$ws.Range("C14").Value = "Coquimbo"
$ws.Range("D14").Value = 44558
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 100112044
$ws.Range("G14").Value = "Perejil"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 2900
$ws.Range("K14").Value = 2300
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2400
$ws.Range("N14").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O14").Value = "Provincia del Elquí"
$ws.Range("P14").Value = 1600
$ws.Range("Q14").Value = 1.5
$ws.Range("R14").Value = "Hortaliza"
